$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12: reuse the border/alignment style already used by rows 5, 6 and 8
# (style index "3") by copying formats from B8:C8, then overwrite the values.
$ws.Range("B8:C8").Copy()
$ws.Range("B12").PasteSpecial(-4122)

$ws.Range("B12").Value = "_C3D-TEMPLATE_2025_FRA (Architecture v0001h)"
$ws.Range("C12").Value = "Modification des styles Plan pour les lignes de coupe de l'édifice (bleu) et les lignes d'élévation de l'édifice (magenta)"

# Scroll the view so row 8 is at the top and select C17, matching the
# author's view state when they saved the workbook.
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C17").Select()
